$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append six new weekly price rows (206-211) for Damasco / Dina at
# Mercado Mayorista Lo Valledor de Santiago, extending the sheet from
# A1:T205 to A1:T211.

# Row 206
$ws.Cells.Item(206,1).Value = 6
$ws.Cells.Item(206,2).Value = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Cells.Item(206,3).Value = 'Metropolitana'
$ws.Cells.Item(206,4).Value = 44911
$ws.Cells.Item(206,4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(206,5).Value = 13
$ws.Cells.Item(206,6).Value = 'Fruta'
$ws.Cells.Item(206,7).Value = 100103
$ws.Cells.Item(206,8).Value = 'Frutos de hueso (carozo)'
$ws.Cells.Item(206,9).Value = 100103003
$ws.Cells.Item(206,10).Value = 'Damasco'
$ws.Cells.Item(206,11).Value = 'Dina'
$ws.Cells.Item(206,12).Value = 'Especial'
$ws.Cells.Item(206,13).Value = 185
$ws.Cells.Item(206,14).Value = 17000
$ws.Cells.Item(206,15).Value = 17000
$ws.Cells.Item(206,16).Value = 17000
$ws.Cells.Item(206,17).Value = '$/caja 16 kilos'
$ws.Cells.Item(206,18).Value = 'Región Metropolitana'
$ws.Cells.Item(206,19).Value = 1062
$ws.Cells.Item(206,20).Value = 16

# Row 207
$ws.Cells.Item(207,1).Value = 6
$ws.Cells.Item(207,2).Value = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Cells.Item(207,3).Value = 'Metropolitana'
$ws.Cells.Item(207,4).Value = 44911
$ws.Cells.Item(207,4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(207,5).Value = 13
$ws.Cells.Item(207,6).Value = 'Fruta'
$ws.Cells.Item(207,7).Value = 100103
$ws.Cells.Item(207,8).Value = 'Frutos de hueso (carozo)'
$ws.Cells.Item(207,9).Value = 100103003
$ws.Cells.Item(207,10).Value = 'Damasco'
$ws.Cells.Item(207,11).Value = 'Dina'
$ws.Cells.Item(207,12).Value = 'Especial'
$ws.Cells.Item(207,13).Value = 80
$ws.Cells.Item(207,14).Value = 17000
$ws.Cells.Item(207,15).Value = 17000
$ws.Cells.Item(207,16).Value = 17000
$ws.Cells.Item(207,17).Value = '$/caja 18 kilos'
$ws.Cells.Item(207,18).Value = 'Región Metropolitana'
$ws.Cells.Item(207,19).Value = 944
$ws.Cells.Item(207,20).Value = 18

# Row 208
$ws.Cells.Item(208,1).Value = 6
$ws.Cells.Item(208,2).Value = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Cells.Item(208,3).Value = 'Metropolitana'
$ws.Cells.Item(208,4).Value = 44911
$ws.Cells.Item(208,4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(208,5).Value = 13
$ws.Cells.Item(208,6).Value = 'Fruta'
$ws.Cells.Item(208,7).Value = 100103
$ws.Cells.Item(208,8).Value = 'Frutos de hueso (carozo)'
$ws.Cells.Item(208,9).Value = 100103003
$ws.Cells.Item(208,10).Value = 'Damasco'
$ws.Cells.Item(208,11).Value = 'Dina'
$ws.Cells.Item(208,12).Value = 'Primera'
$ws.Cells.Item(208,13).Value = 300
$ws.Cells.Item(208,14).Value = 14000
$ws.Cells.Item(208,15).Value = 15000
$ws.Cells.Item(208,16).Value = 14500
$ws.Cells.Item(208,17).Value = '$/caja 16 kilos'
$ws.Cells.Item(208,18).Value = 'Región Metropolitana'
$ws.Cells.Item(208,19).Value = 906
$ws.Cells.Item(208,20).Value = 16

# Row 209
$ws.Cells.Item(209,1).Value = 6
$ws.Cells.Item(209,2).Value = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Cells.Item(209,3).Value = 'Metropolitana'
$ws.Cells.Item(209,4).Value = 44911
$ws.Cells.Item(209,4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(209,5).Value = 13
$ws.Cells.Item(209,6).Value = 'Fruta'
$ws.Cells.Item(209,7).Value = 100103
$ws.Cells.Item(209,8).Value = 'Frutos de hueso (carozo)'
$ws.Cells.Item(209,9).Value = 100103003
$ws.Cells.Item(209,10).Value = 'Damasco'
$ws.Cells.Item(209,11).Value = 'Dina'
$ws.Cells.Item(209,12).Value = 'Primera'
$ws.Cells.Item(209,13).Value = 100
$ws.Cells.Item(209,14).Value = 15000
$ws.Cells.Item(209,15).Value = 15000
$ws.Cells.Item(209,16).Value = 15000
$ws.Cells.Item(209,17).Value = '$/caja 18 kilos'
$ws.Cells.Item(209,18).Value = 'Región Metropolitana'
$ws.Cells.Item(209,19).Value = 833
$ws.Cells.Item(209,20).Value = 18

# Row 210
$ws.Cells.Item(210,1).Value = 6
$ws.Cells.Item(210,2).Value = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Cells.Item(210,3).Value = 'Metropolitana'
$ws.Cells.Item(210,4).Value = 44911
$ws.Cells.Item(210,4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(210,5).Value = 13
$ws.Cells.Item(210,6).Value = 'Fruta'
$ws.Cells.Item(210,7).Value = 100103
$ws.Cells.Item(210,8).Value = 'Frutos de hueso (carozo)'
$ws.Cells.Item(210,9).Value = 100103003
$ws.Cells.Item(210,10).Value = 'Damasco'
$ws.Cells.Item(210,11).Value = 'Dina'
$ws.Cells.Item(210,12).Value = 'Segunda'
$ws.Cells.Item(210,13).Value = 250
$ws.Cells.Item(210,14).Value = 10000
$ws.Cells.Item(210,15).Value = 10000
$ws.Cells.Item(210,16).Value = 10000
$ws.Cells.Item(210,17).Value = '$/caja 16 kilos'
$ws.Cells.Item(210,18).Value = 'Región Metropolitana'
$ws.Cells.Item(210,19).Value = 625
$ws.Cells.Item(210,20).Value = 16

# Row 211
$ws.Cells.Item(211,1).Value = 6
$ws.Cells.Item(211,2).Value = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Cells.Item(211,3).Value = 'Metropolitana'
$ws.Cells.Item(211,4).Value = 44911
$ws.Cells.Item(211,4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(211,5).Value = 13
$ws.Cells.Item(211,6).Value = 'Fruta'
$ws.Cells.Item(211,7).Value = 100103
$ws.Cells.Item(211,8).Value = 'Frutos de hueso (carozo)'
$ws.Cells.Item(211,9).Value = 100103003
$ws.Cells.Item(211,10).Value = 'Damasco'
$ws.Cells.Item(211,11).Value = 'Dina'
$ws.Cells.Item(211,12).Value = 'Segunda'
$ws.Cells.Item(211,13).Value = 150
$ws.Cells.Item(211,14).Value = 11000
$ws.Cells.Item(211,15).Value = 11000
$ws.Cells.Item(211,16).Value = 11000
$ws.Cells.Item(211,17).Value = '$/caja 18 kilos'
$ws.Cells.Item(211,18).Value = 'Región Metropolitana'
$ws.Cells.Item(211,19).Value = 611
$ws.Cells.Item(211,20).Value = 18
